$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: move the "_GoBack" bookmark from between the "t" / "utorial"
# runs in the tutorial paragraph to the end of the following (empty)
# paragraph.
# ------------------------------------------------------------------
$tutorialRange = $d.Content
$tutorialRange.Find.Execute("ecoED Champions Training 2018", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$targetPara = $null
$tutorialParaSeen = $false
foreach ($p in $d.Paragraphs) {
    if ($tutorialParaSeen) {
        $targetPara = $p
        break
    }
    if (($p.Range.Start -le $tutorialRange.Start) -and ($p.Range.End -ge $tutorialRange.End)) {
        $tutorialParaSeen = $true
    }
}

$bmRange = $targetPara.Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# Change 2: split the "install_github("GregGuerin/ausplotsR")" run into
# three runs -- "install_github(", "ternaustralia" and "/ausplotsR")" --
# replacing "GregGuerin" with "ternaustralia" in the middle.
# ------------------------------------------------------------------
$matchRange = $d.Content
$matchRange.Find.Execute("GregGuerin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Temporarily toggle bold on the matched range so the engine keeps it as
# its own run once the text is replaced, then clear the bold again so the
# final formatting matches the surrounding runs exactly.
$matchRange.Font.Bold = 1
$matchRange.Text = "ternaustralia"

$newStart = $matchRange.Start
$newRange = $d.Range($newStart, $newStart + 13)
$newRange.Font.Bold = 0
